$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 20 updates - mirror style/structure of row 19 (the sibling "done" row)
$ws.Range("B19:E19").Copy()
$ws.Range("B20:E20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("G19:H19").Copy()
$ws.Range("G20:H20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("F19").Copy()
$ws.Range("F20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("E20").Value = 15
$ws.Range("F20").Value = Get-Date -Year 2021 -Month 10 -Day 26 -Hour 0 -Minute 0 -Second 0
$ws.Range("H20").Value = "need to be commented"

$excel.CutCopyMode = $false

# Update the current selection to reflect where the user left off
$ws.Range("G20").Select()
